# Backup QR Scanner data - 4/6/2025, 4:53:45 PM
# Append the newest QR-scan log entries (rows 9-15) to the "Jzbdhd" sheet,
# mirroring the same layout/typing already used by the existing rows
# (Number = numeric, Student ID / Location / Log Date / Log Time = text).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Jzbdhd")

$newEntries = @(
    @{ Number = 8;  StudentId = "231249"; Location = "Jzbdhd"; LogDate = "2025-04-06"; LogTime = "16:53:19" },
    @{ Number = 9;  StudentId = "231249"; Location = "Jzbdhd"; LogDate = "2025-04-06"; LogTime = "16:53:22" },
    @{ Number = 10; StudentId = "231249"; Location = "Jzbdhd"; LogDate = "2025-04-06"; LogTime = "16:53:25" },
    @{ Number = 11; StudentId = "231249"; Location = "Jzbdhd"; LogDate = "2025-04-06"; LogTime = "16:53:30" },
    @{ Number = 12; StudentId = "231249"; Location = "Jzbdhd"; LogDate = "2025-04-06"; LogTime = "16:53:33" },
    @{ Number = 13; StudentId = "231249"; Location = "Jzbdhd"; LogDate = "2025-04-06"; LogTime = "16:53:36" },
    @{ Number = 14; StudentId = "231249"; Location = "Jzbdhd"; LogDate = "2025-04-06"; LogTime = "16:53:40" }
)

$startRow = 9
$row = $startRow
foreach ($entry in $newEntries) {
    # Column A: plain number
    $ws.Cells.Item($row, 1).Value = $entry.Number

    # Columns B, D, E look numeric/date-like (e.g. "231249", "16:53:19") -
    # prefix with an apostrophe so Excel keeps them as text, just like the
    # rest of the sheet's Student ID / Log Date / Log Time columns.
    $ws.Cells.Item($row, 2).Value = "'" + $entry.StudentId
    $ws.Cells.Item($row, 3).Value = $entry.Location
    $ws.Cells.Item($row, 4).Value = "'" + $entry.LogDate
    $ws.Cells.Item($row, 5).Value = "'" + $entry.LogTime

    $row++
}
